$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'330.76"
$ws.Range("E2").Value = "'0.47%"
$ws.Range("G2").Value = "'21"
$ws.Range("D3").Value = "'44.53"
$ws.Range("E3").Value = "'1.03%"
$ws.Range("G3").Value = "'21"
$ws.Range("D4").Value = "'5.535"
$ws.Range("G4").Value = "'21"
$ws.Range("E5").Value = "'0.94%"
$ws.Range("G5").Value = "'21"
$ws.Range("D6").Value = "'2.056"
$ws.Range("E6").Value = "'3.97%"
$ws.Range("G6").Value = "'21"
$ws.Range("D7").Value = "'0.9752"
$ws.Range("E7").Value = "'2.42%"
$ws.Range("G7").Value = "'21"
$ws.Range("D8").Value = "'0.1117"
$ws.Range("E8").Value = "'-4.13%"
$ws.Range("G8").Value = "'21"
$ws.Range("D9").Value = "'0.1898"
$ws.Range("E9").Value = "'2.25%"
$ws.Range("G9").Value = "'21"
$ws.Range("D10").Value = "'10.24"
$ws.Range("E10").Value = "'-13.80%"
$ws.Range("G10").Value = "'21"
$ws.Range("D11").Value = "'0.1002"
$ws.Range("E11").Value = "'2.19%"
$ws.Range("G11").Value = "'21"
$ws.Range("D12").Value = "'0.04736"
$ws.Range("E12").Value = "'0.04%"
$ws.Range("G12").Value = "'21"
$ws.Range("D13").Value = "'0.1057"
$ws.Range("E13").Value = "'-0.98%"
$ws.Range("G13").Value = "'21"
$ws.Range("D14").Value = "'0.001260"
$ws.Range("E14").Value = "'-1.91%"
$ws.Range("G14").Value = "'21"
$ws.Range("D15").Value = "'0.04099"
$ws.Range("E15").Value = "'-3.24%"
$ws.Range("G15").Value = "'21"
$ws.Range("D16").Value = "'0.006011"
$ws.Range("E16").Value = "'1.55%"
$ws.Range("G16").Value = "'21"
$ws.Range("D17").Value = "'3.347"
$ws.Range("G17").Value = "'21"
$ws.Range("D18").Value = "'4.428"
$ws.Range("E18").Value = "'2.24%"
$ws.Range("G18").Value = "'21"
$ws.Range("E19").Value = "'2.94%"
$ws.Range("G19").Value = "'21"
$ws.Range("D20").Value = "'0.3351"
$ws.Range("E20").Value = "'-3.57%"
$ws.Range("G20").Value = "'21"
$ws.Range("D21").Value = "'0.1389"
$ws.Range("E21").Value = "'-1.38%"
$ws.Range("G21").Value = "'21"
$ws.Range("D22").Value = "'0.2566"
$ws.Range("E22").Value = "'2.33%"
$ws.Range("G22").Value = "'21"
$ws.Range("D23").Value = "'0.001303"
$ws.Range("E23").Value = "'3.87%"
$ws.Range("G23").Value = "'21"
$ws.Range("D24").Value = "'0.004399"
$ws.Range("E24").Value = "'1.84%"
$ws.Range("G24").Value = "'21"
$ws.Range("D25").Value = "'0.0001276"
$ws.Range("E25").Value = "'7.18%"
$ws.Range("G25").Value = "'21"
$ws.Range("D26").Value = "'0.0003730"
$ws.Range("E26").Value = "'-6.26%"
$ws.Range("G26").Value = "'21"
$ws.Range("G27").Value = "'21"
$ws.Range("G28").Value = "'21"
$ws.Range("G29").Value = "'21"
$ws.Range("G30").Value = "'21"
$ws.Range("G31").Value = "'21"
$ws.Range("G32").Value = "'21"
$ws.Range("G33").Value = "'21"
$ws.Range("G34").Value = "'21"
$ws.Range("G35").Value = "'21"
$ws.Range("G36").Value = "'21"
$ws.Range("G37").Value = "'21"
$ws.Range("D38").Value = "'0.02687"
$ws.Range("E38").Value = "'1.04%"
$ws.Range("G38").Value = "'21"
$ws.Range("D39").Value = "'0.05679"
$ws.Range("E39").Value = "'2.30%"
$ws.Range("G39").Value = "'21"
$ws.Range("D40").Value = "'0.007621"
$ws.Range("E40").Value = "'0.74%"
$ws.Range("G40").Value = "'21"
$ws.Range("D41").Value = "'0.1422"
$ws.Range("E41").Value = "'1.00%"
$ws.Range("G41").Value = "'21"
$ws.Range("D42").Value = "'0.007500"
$ws.Range("E42").Value = "'-7.19%"
$ws.Range("G42").Value = "'21"
$ws.Range("E43").Value = "'-3.16%"
$ws.Range("G43").Value = "'21"
$ws.Range("D44").Value = "'0.008287"
$ws.Range("E44").Value = "'-6.84%"
$ws.Range("G44").Value = "'21"
$ws.Range("D45").Value = "'0.00007043"
$ws.Range("E45").Value = "'-2.90%"
$ws.Range("G45").Value = "'21"
$ws.Range("E46").Value = "'-0.30%"
$ws.Range("G46").Value = "'21"
$ws.Range("D47").Value = "'0.0005784"
$ws.Range("E47").Value = "'-0.47%"
$ws.Range("G47").Value = "'21"
$ws.Range("D48").Value = "'0.002513"
$ws.Range("E48").Value = "'10.67%"
$ws.Range("G48").Value = "'21"
$ws.Range("D49").Value = "'0.003531"
$ws.Range("E49").Value = "'-25.27%"
$ws.Range("G49").Value = "'21"
$ws.Range("E50").Value = "'-0.30%"
$ws.Range("G50").Value = "'21"
$ws.Range("E51").Value = "'-0.30%"
$ws.Range("G51").Value = "'21"
